$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "42.501.36"
$ws.Range("E2").Value = "  -0.46%  "
$ws.Range("D3").Value = "2.287.08"
$ws.Range("E3").Value = "  -0.51%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "0.999"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  -0.38%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "311.05"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.86%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "103.54"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -0.96%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.622"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.33%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("E9").Value = "  -1.71%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "38.73"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -3.69%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0900"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.01%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "8.20"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -3.85%  "
$ws.Range("E13").Value = "  +0.56%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.969"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.39%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.22"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.43%  "
$ws.Range("D16").Value = "2.632.42"
$ws.Range("E16").Value = "  -0.65%  "
$ws.Range("D17").Value = "2.291.34"
$ws.Range("E17").Value = "  -0.17%  "
$ws.Range("D18").Value = "42.464.82"
$ws.Range("E18").Value = "  -0.51%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.29"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.40%  "
$ws.Range("E20").Value = "  -1.52%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "13.34"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "73.20"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.22%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "269.89"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.16%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "3.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.58%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.17"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -3.32%  "
$ws.Range("E26").Value = "  +0.28%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.73"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.81%  "
$ws.Range("B28").Value = "Filecoin"
$ws.Range("C28").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "6.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +12.50%  "
$ws.Range("B29").Value = "Toncoin"
$ws.Range("C29").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.28"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -2.28%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "22.30"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -1.31%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "35.87"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -7.28%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "164.33"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0846"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -4.12%  "
$ws.Range("E34").Value = "  -2.95%  "
$ws.Range("E35").Value = "  +0.33%  "
$ws.Range("E36").Value = "  -3.18%  "
$ws.Range("E37").Value = "  -3.26%  "
$ws.Range("E38").Value = "  -3.50%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.75"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.01%  "
$ws.Range("B40").Value = "NEARProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.61"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "
$ws.Range("B41").Value = "BitcoinSV"
$ws.Range("C41").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "111.76"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +19.32%  "
$ws.Range("E42").Value = "  +0.78%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "70.02"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -0.57%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.225"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.66%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "12.07"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -2.58%  "
$ws.Range("D47").Value = "1.720.55"
$ws.Range("E47").Value = "  +7.69%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "110.10"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -3.47%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "76.94"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -5.75%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "5.16"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -2.38%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "8.64"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -3.38%  "
